$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the (erroneous) leading column A; this shifts B:F left to A:E,
# dropping the old column A values/style entirely.
$ws.Columns("A").Delete()

# Fix the header text: "MODEL_CONDITION" -> "MODELCONDITION".
# After the shift, that header now lives in column D.
$ws.Range("D1").Value = "MODELCONDITION"
